$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before column E (old E..Z shift right to I..AD)
$ws.Range("E1:H1").EntireColumn.Insert()

# New header row values (row 1)
$ws.Range("E1").Value2 = "URL_ICONS"
$ws.Range("F1").Value2 = "PASSWORD_ICONS"
$ws.Range("G1").Value2 = "KODE_CABANG"
$ws.Range("H1").Value2 = "NOMOR_TERMINAL"

# New data row values (row 2)
$ws.Range("E2").Value2 = "http://192.168.150.186/alternity/"
$ws.Range("F2").Value2 = 1
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value2 = "089"
$ws.Range("H2").Value2 = 259

# Style the new E2/G2 cells like the D2 hyperlink-derived style, but plain font + left/center alignment
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("G2").Style = "Hyperlink"
$ws.Range("E2:H2").HorizontalAlignment = -4131
$ws.Range("E2:H2").VerticalAlignment = -4108

# Column widths for the newly inserted columns
$ws.Columns("E").ColumnWidth = 29.88
$ws.Columns("F").ColumnWidth = 17.17
$ws.Columns("G").ColumnWidth = 13.74
$ws.Columns("H").ColumnWidth = 17.17

# New hyperlink on E2
$ws.Hyperlinks.Add($ws.Range("E2"), "http://192.168.150.186/alternity/")

# Reset the top-left scroll position back to A1 (was D1)
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1

Write-Output "done"
